$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2020" data column (X) is appended right after the existing "2019"
# column (W). Each new cell should carry the same number format/font/border
# as the corresponding cell in column W for that row, so copy the format
# over first and then set the new value.

$years = @{
    4  = 2020
    5  = 45.3
    6  = 48.2
    7  = 43.6
    8  = 48.8
    9  = 41.5
    10 = 49.7
    11 = 46.7
    12 = 36.5
    13 = 29.6
    14 = 54.7
    15 = 51.6
    16 = 47.2
}

foreach ($row in 4..16) {
    $ws.Range("W$row").Copy()
    $ws.Range("X$row").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("X$row").Value = $years[$row]
}

$excel.CutCopyMode = $false

# Restore the (non-data) active-cell selection recorded in the saved view.
[void]$ws.Range("AI21").Select()
